$d = $word.ActiveDocument

# 1. "From author annotations" -> "For author annotations"
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r") -eq "From author annotations") {
        $r = $p.Range
        $r.End = $r.End - 1
        $r.Text = "For author annotations"
    }
}

# 2. Insert two new list paragraphs ("Credit Title", "Credit Authors")
#    right after the "Resource Title" paragraph, matching its list formatting
#    (ilvl=1, numId=1004).
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r") -eq "Resource Title") {
        $p.Range.InsertParagraphAfter()
        $creditTitle = $p.Next()
        $creditTitle.Range.Text = "Credit Title"

        $creditTitle.Range.InsertParagraphAfter()
        $creditAuthors = $creditTitle.Next()
        $creditAuthors.Range.Text = "Credit Authors"
    }
}
